$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row changes
$ws.Range("H1").Value = "Salt"
$ws.Range("N1").Value = "Card Art"
$ws.Range("O1").Value = "Picture Location"

# Row 2 - Tomato
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("N2").Value = "Tomato.png"
$ws.Range("O2").Value = "ING001"

# Row 3 - Heirloom Tomato
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 3
$ws.Range("N3").Value = "HeirloomTomato.jpg"
$ws.Range("O3").Value = "ING002"

# Row 4 - Spaghetti
$ws.Range("N4").Value = "Spaghetti.jpg"
$ws.Range("O4").Value = "ING003"

# Row 5 - Mom's Spaghetti
$ws.Range("N5").Value = "Spaghetti.jpg"
$ws.Range("O5").Value = "ING004"

# Column width for new column N
$ws.Columns.Item(14).ColumnWidth = 18

# Update selection to match target state
[void]$ws.Range("C41").Select()
